$wb = $excel.ActiveWorkbook
$wsContacts = $wb.Worksheets.Item("ContactsTestData")
$wsOrg = $wb.Worksheets.Item("OrganizationTestData")

# --- Add new reference/list data to OrganizationTestData (columns E-M, rows 7-9) ---
# Values are entered in this specific order so that the shared-string table
# ends up with the same ordering as the target workbook.
$wsOrg.Range("F7").Value = "Banking"
$wsOrg.Range("E7").Value = "Chemicals"
$wsOrg.Range("G7").Value = "Consulting"

$wsOrg.Range("E8").Value = "Analyst"
$wsOrg.Range("F8").Value = "Competiitor"
$wsOrg.Range("G8").Value = "Customer"
$wsOrg.Range("H8").Value = "Integrator"
$wsOrg.Range("I8").Value = "Partner"
$wsOrg.Range("J8").Value = "Press"
$wsOrg.Range("K8").Value = "Prospect"
$wsOrg.Range("L8").Value = "Reseller"
$wsOrg.Range("M8").Value = "Other"

$wsOrg.Range("E9").Value = "Active"
$wsOrg.Range("F9").Value = "Market Failed"
$wsOrg.Range("G9").Value = "Project Cancelled"
$wsOrg.Range("H9").Value = "Shutdown"

# --- Set the width of the new columns F and G ---
$wsOrg.Columns.Item(6).ColumnWidth = 12.42
$wsOrg.Columns.Item(7).ColumnWidth = 18.59

# --- Update sheet selections / active sheet ---
$wsContacts.Activate()
$wsContacts.Range("F23").Select()

$wsOrg.Activate()
$wsOrg.Range("F23").Select()
